{"js": "// Update the copyright year on the cover page from 2021 to 2022.\nconst results = context.document.body.search(\"2021\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"2022\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the copyright year on the cover page from 2021 to 2022.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"2021\"\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $find.Parent.Text = \"2022\"\n}\n"}
